# Append the "16 December" journal entry to the end of the Wiki Procesboek.
# The new block consists of: two blank spacer paragraphs, a bold/16pt date
# heading ("16 December"), another blank spacer paragraph, and the body
# paragraph describing the day's bug fixes.
#
# We build the block as a single raw-OOXML fragment and insert it with
# Range.InsertXML at the very end of the document's story. Doing it as one
# multi-paragraph fragment (rather than several separate
# InsertParagraphAfter calls) makes every paragraph break "real" -
# including the blank paragraphs, which end up with no run at all, just
# like the blank paragraphs already used throughout this document.

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$insertionPoint = $lastPara.Range.End
$target = $d.Range($insertionPoint, $insertionPoint)

$newBlockXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="nl-NL"/></w:rPr><w:t>16 December</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Ik ben net klaar met Commerce dus dacht ik laat ik nog even wiki controleren. Ik kwam erachter dat ik voor de entry pagina’s niet wiki/entry had maar alleen de entry naam achter de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">, dit heb ik nog aangepast. Ook als ik een pagina probeerde aan te maken die al </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">bestond, schreef hij hem over. In de opdracht stond dat als de gebruiker een pagina probeerde aan te maken die al bestond, degene naar een error page gestuurd zou moeten worden. Dit heb ik ook nog toegevoegd. Nu zijn alle functionaliteiten zoals de opdracht voorschrijft. </w:t></w:r></w:p>'

$target.InsertXML($newBlockXml)

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
